$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("aquarius-plus-top-pos")
$ws.Activate()

# --- Component placement edits (Designator D3, J1-J4, J9, JP1-JP3 moved) ---
# D3  (row 52): Mid X changed
$ws.Range("B52").Value = 76.77

# J1  (row 53): Mid Y changed
$ws.Range("C53").Value = -45.8

# J2  (row 54): Mid X, Mid Y, Rotation changed
$ws.Range("B54").Value = 97.45
$ws.Range("C54").Value = -52.5
$ws.Range("E54").Value = 180

# J3  (row 55): Mid Y changed
$ws.Range("C55").Value = -48.65

# J4  (row 56): Mid X, Mid Y changed
$ws.Range("B56").Value = 149.19999999999999
$ws.Range("C56").Value = -54.5

# J9  (row 57): Mid Y changed
$ws.Range("C57").Value = -109.02

# JP1 (row 60): Mid Y changed
$ws.Range("C60").Value = -55.25

# JP2 (row 61): Mid Y changed
$ws.Range("C61").Value = -55.25

# JP3 (row 62): Mid X changed
$ws.Range("B62").Value = 86.05

# --- View state: scroll/zoom/selection to match where the author was working ---
$window = $excel.ActiveWindow
$window.ScrollRow = 43
$window.ScrollColumn = 1
$window.Zoom = 175
$ws.Range("C57").Select()
